$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 218; this shifts the existing rows
# 218-225 down to 219-226 and expands the used range to A1:R226.
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new record.
$ws.Cells.Item(218, 1).Value = 5
$ws.Cells.Item(218, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(218, 3).Value = "Maule"
$ws.Cells.Item(218, 4).Value = 44509
$ws.Cells.Item(218, 5).Value = 7
$ws.Cells.Item(218, 6).Value = 100112032
$ws.Cells.Item(218, 7).Value = "Zapallo italiano"
$ws.Cells.Item(218, 8).Value = "Sin especificar"
$ws.Cells.Item(218, 9).Value = "Primera"
$ws.Cells.Item(218, 10).Value = 500
$ws.Cells.Item(218, 11).Value = 7000
$ws.Cells.Item(218, 12).Value = 7000
$ws.Cells.Item(218, 13).Value = 7000
$ws.Cells.Item(218, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(218, 15).Value = "Región del Maule"
$ws.Cells.Item(218, 16).Value = 117
$ws.Cells.Item(218, 17).Value = 60
$ws.Cells.Item(218, 18).Value = "Hortaliza"
